$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 19, shifting existing rows 19-127 down to 20-128.
$ws.Rows.Item(19).Insert()

# Populate the newly inserted row 19 with the new weekly entry.
$ws.Cells.Item(19, 1).Value  = 11
$ws.Cells.Item(19, 2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(19, 3).Value  = "Bíobío"
$ws.Cells.Item(19, 4).Value  = 44561
$ws.Cells.Item(19, 5).Value  = 8
$ws.Cells.Item(19, 6).Value  = "Fruta"
$ws.Cells.Item(19, 7).Value  = 100108
$ws.Cells.Item(19, 8).Value  = "Tropicales y subtropicales"
$ws.Cells.Item(19, 9).Value  = 100108005
$ws.Cells.Item(19, 10).Value = "Piña"
$ws.Cells.Item(19, 11).Value = "Caramelo"
$ws.Cells.Item(19, 12).Value = "Primera"
$ws.Cells.Item(19, 13).Value = 200
$ws.Cells.Item(19, 14).Value = 15000
$ws.Cells.Item(19, 15).Value = 16000
$ws.Cells.Item(19, 16).Value = 15500
$ws.Cells.Item(19, 17).Value = "$/caja 12 unidades"
$ws.Cells.Item(19, 18).Value = "Ecuador"
$ws.Cells.Item(19, 19).Value = 1292
$ws.Cells.Item(19, 20).Value = 12
